$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename the second sheet ---
$ws2.Name = "MasterDataCreation"

# --- sheet1 (MasterTestDataSheet) row 2 data updates ---
# Keep the existing quotePrefix style (s=3) by prefixing values with a
# leading apostrophe, which Excel treats as a literal-text marker rather
# than part of the stored value.
$ws1.Range("B2").Value = "'MasterDataCreation"
$ws1.Range("D2").Value = "'..\\JunoAutomation\src\resources\\Juno_TestDataSheet.xlsx"

# --- sheet2 (MasterDataCreation) new columns D:G ---
# Copy header style (B1) onto the new header cells, and the data style
# (B2/C2) onto the new JAMS cell, before assigning their values so the
# copied formatting isn't clobbered by the value assignment.
$ws2.Range("B1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("E1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$ws2.Range("C2").Copy()
$ws2.Range("D2").PasteSpecial(-4122)

$ws2.Range("D1").Value = "App"
$ws2.Range("E1").Value = "NonPrescriber"
$ws2.Range("F1").Value = "Organisation"
$ws2.Range("G1").Value = "Site"

# D2 inherited the quotePrefix style (s=6) from C2 via PasteSpecial; a plain
# Value assignment would silently drop that quotePrefix flag (same trick as
# above), so prefix with an apostrophe to keep the style intact.
$ws2.Range("D2").Value = "'JAMS"
$ws2.Range("E2").Value = "Auto_TueJan081543172019"
$ws2.Range("F2").Value = "Auto_TueJan081543172019"
$ws2.Range("G2").Value = "Auto_TueJan0815431720975"

# --- Column widths for the new sheet2 columns (E/F ~23.4 chars, G ~24.4) ---
$ws2.Columns.Item(5).ColumnWidth = 22.5
$ws2.Columns.Item(6).ColumnWidth = 22.5
$ws2.Columns.Item(7).ColumnWidth = 23.5

# --- Selections / active sheet ---
$ws1.Range("B3").Select()
$ws2.Range("F10").Select()
$ws2.Activate()

Write-Output "done"
